# Update the "dSF" (column F) values for several rows as part of a
# repull/push of the underlying data and mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    6  = 9
    7  = 2
    8  = -2
    11 = 1
    18 = 1
    26 = -3
    35 = -2
    40 = -1
    41 = -3
    54 = 1
    57 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
